$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / "want to go" count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12895
$ws1.Range("F3").Value = 636
$ws1.Range("F5").Value = 44
$ws1.Range("F6").Value = 331
$ws1.Range("F8").Value = 237
$ws1.Range("F9").Value = 13041
$ws1.Range("F11").Value = 36
$ws1.Range("F12").Value = 5292
$ws1.Range("F13").Value = 552
$ws1.Range("F16").Value = 35
$ws1.Range("F17").Value = 1203
$ws1.Range("F21").Value = 2865
$ws1.Range("F22").Value = 6221
$ws1.Range("F23").Value = 1167
$ws1.Range("F24").Value = 3641
$ws1.Range("F26").Value = 49

# Sheet "全部类型" (All Types) - update column F (想去人数 / "want to go" count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12895
$ws4.Range("F3").Value = 636
$ws4.Range("F5").Value = 44
$ws4.Range("F6").Value = 331
$ws4.Range("F9").Value = 237
$ws4.Range("F10").Value = 13041
$ws4.Range("F12").Value = 36
$ws4.Range("F13").Value = 5292
$ws4.Range("F14").Value = 552
$ws4.Range("F17").Value = 35
$ws4.Range("F18").Value = 1203
$ws4.Range("F22").Value = 2865
$ws4.Range("F24").Value = 6221
$ws4.Range("F25").Value = 1167
$ws4.Range("F26").Value = 3641
$ws4.Range("F28").Value = 49
